$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicated header row that was accidentally appended at the
# bottom of the roster (row 43 mirrored row 2's header labels).
$ws.Rows(43).Delete()

# Widen column A so the longest player name ("Red Schoendienst HOF") is
# fully visible.
$ws.Columns("A").ColumnWidth = 27.6667

# Leave the view scrolled down with the (now blank) row 43 selected, as if
# the user had just deleted that row and was about to start building a new
# scoresheet form below the roster.
$ws.Rows(43).Select()
